# Add team win/loss/tie record columns (AD:AF) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should carry the same bold/bordered/centered header
# style already used across row 1 (style index 1, e.g. the "Unnamed: 28"
# cell in AC1). Copy that formatting onto the three new header cells,
# then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-57) gets the same team record: 73 wins, 89 losses,
# 0 ties.
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
